$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the existing data row (row 2) down into rows 3-6,
# then fill each new row with the same values as row 2 (import of
# several UO rows at once).
$ws.Range("A2:S2").Copy()

for ($r = 3; $r -le 6; $r++) {
    $target = $ws.Range("A" + $r + ":S" + $r)
    $target.PasteSpecial(-4122)

    for ($col = 1; $col -le 19; $col++) {
        $ws.Cells.Item($r, $col).Value2 = $ws.Cells.Item(2, $col).Value2
    }
}

$ws.Range("M13").Select()
